# Insert a new data row before current row 32, shifting existing rows 32:159 down to 33:160,
# and populate the new row 32 with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 32 (pushes old row 32.. down by one)
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with the new record's values
$ws.Cells.Item(32, 1).Value = 5
$ws.Cells.Item(32, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(32, 3).Value = "Maule"
$ws.Cells.Item(32, 4).Value = (Get-Date -Year 2022 -Month 4 -Day 25 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(32, 5).Value = 7
$ws.Cells.Item(32, 6).Value = 100112017
$ws.Cells.Item(32, 7).Value = "Apio"
$ws.Cells.Item(32, 8).Value = "Americana (o)"
$ws.Cells.Item(32, 9).Value = "Primera"
$ws.Cells.Item(32, 10).Value = 600
$ws.Cells.Item(32, 11).Value = 7000
$ws.Cells.Item(32, 12).Value = 7000
$ws.Cells.Item(32, 13).Value = 7000
$ws.Cells.Item(32, 14).Value = "`$/docena de matas"
$ws.Cells.Item(32, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(32, 16).Value = 1167
$ws.Cells.Item(32, 17).Value = 6
$ws.Cells.Item(32, 18).Value = "Hortaliza"
